$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D:E to text so numeric-looking strings keep their exact formatting
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.486.08"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "1.934.65"
$ws.Range("E3").Value = "  +4.42%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "240.72"
$ws.Range("E5").Value = "  +3.12%  "
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "0.4747"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "0.2875"
$ws.Range("E8").Value = "  +4.31%  "
$ws.Range("D9").Value = "0.06640"
$ws.Range("E9").Value = "  +4.67%  "
$ws.Range("D10").Value = "19.26"
$ws.Range("E10").Value = "  +7.55%  "
$ws.Range("D11").Value = "106.71"
$ws.Range("E11").Value = "  +25.91%  "
$ws.Range("D12").Value = "1.926.59"
$ws.Range("E12").Value = "  -8.01%  "
$ws.Range("D13").Value = "0.07628"
$ws.Range("E13").Value = "  +2.18%  "
$ws.Range("E14").Value = "  +3.91%  "
$ws.Range("D15").Value = "0.6648"
$ws.Range("E15").Value = "  +6.38%  "
$ws.Range("D16").Value = "309.49"
$ws.Range("E16").Value = "  +26.05%  "
$ws.Range("D17").Value = "30.490.04"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "13.04"
$ws.Range("E18").Value = "  +2.76%  "
$ws.Range("D19").Value = "0.9998"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "0.000007572"
$ws.Range("E20").Value = "  +3.16%  "
$ws.Range("D21").Value = "2.175.15"
$ws.Range("E21").Value = "  +2.78%  "
$ws.Range("D22").Value = "5.310"
$ws.Range("E22").Value = "  +7.78%  "
$ws.Range("D23").Value = "0.9994"
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").Value = "6.311"
$ws.Range("E24").Value = "  +6.62%  "
$ws.Range("D25").Value = "9.320"
$ws.Range("E25").Value = "  +2.58%  "
$ws.Range("D26").Value = "167.71"
$ws.Range("E26").Value = "  +2.00%  "
$ws.Range("D27").Value = "20.34"
$ws.Range("E27").Value = "  +13.15%  "
$ws.Range("D28").Value = "2.063"
$ws.Range("E28").Value = "  +9.85%  "
$ws.Range("D29").Value = "0.1109"
$ws.Range("E29").Value = "  +8.00%  "
$ws.Range("D30").Value = "1.368"
$ws.Range("E30").Value = "  +1.49%  "
$ws.Range("D31").Value = "4.120"
$ws.Range("E31").Value = "  +1.83%  "
$ws.Range("E32").Value = "  +2.71%  "
$ws.Range("D33").Value = "0.05044"
$ws.Range("E33").Value = "  +4.33%  "
$ws.Range("D34").Value = "0.7426"
$ws.Range("E34").Value = "  +6.25%  "
$ws.Range("D35").Value = "1.156"
$ws.Range("E35").Value = "  +2.28%  "
$ws.Range("D36").Value = "2.762"
$ws.Range("E36").Value = "  +2.32%  "
$ws.Range("E37").Value = "  +3.38%  "
$ws.Range("D38").Value = "2.691"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("D39").Value = "2.050"
$ws.Range("E39").Value = "  +2.79%  "
$ws.Range("D40").Value = "0.8800"
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("D41").Value = "107.66"
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("D42").Value = "70.61"
$ws.Range("E42").Value = "  +11.60%  "
$ws.Range("D43").Value = "5.810"
$ws.Range("E43").Value = "  +5.40%  "
$ws.Range("D44").Value = "0.9993"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").Value = "0.4191"
$ws.Range("E45").Value = "  +3.12%  "
$ws.Range("D46").Value = "7.301"
$ws.Range("E46").Value = "  +1.77%  "
$ws.Range("D47").Value = "9.264"
$ws.Range("E47").Value = "  +8.35%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "34.93"
$ws.Range("E48").Value = "  +2.63%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "0.1212"
$ws.Range("E49").Value = "  +0.85%  "
$ws.Range("D50").Value = "0.05624"
$ws.Range("E50").Value = "  +2.16%  "
$ws.Range("E51").Value = "  +4.47%  "

# Restore default (unstyled) cell style so XML matches original (no explicit s attr)
$ws.Range("D2:E51").Style = "Normal"
